# Estadisticos Segundo Parcial 26 Mayo
#
# The "Rescatables" sheet previously listed a single rescatable student
# (SOTO ZOPILLAXTLE, LUZ ARIANA). This update keeps that student (now on
# row 3) and adds two more rescatable students: one on row 2
# (OSORIO HERNANDEZ, AYLIN ABIGAIL) and a new one on row 4
# (CORTES MENDEZ, ELIZABETH MADAI).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# --- Column A (NC / control numbers) ---
$ws.Cells.Item(2, 1).Value = 24330051920187
$ws.Cells.Item(3, 1).Value = 24330051920201
$ws.Cells.Item(4, 1).Value = 24330051920398

# --- Column B (Paterno) ---
$ws.Cells.Item(2, 2).Value = "OSORIO"
$ws.Cells.Item(4, 2).Value = "CORTES"
$ws.Cells.Item(3, 2).Value = "SOTO"

# --- Column C (Materno) ---
$ws.Cells.Item(2, 3).Value = "HERNANDEZ"
$ws.Cells.Item(4, 3).Value = "MENDEZ"
$ws.Cells.Item(3, 3).Value = "ZOPILLAXTLE"

# --- Column D (Nombres) ---
$ws.Cells.Item(2, 4).Value = "AYLIN ABIGAIL"
$ws.Cells.Item(4, 4).Value = "ELIZABETH MADAI"
$ws.Cells.Item(3, 4).Value = "LUZ ARIANA"

# --- Column E (Nombre_Largo / materia) ---
$ws.Cells.Item(2, 5).Value = "Cultura digital II"
$ws.Cells.Item(3, 5).Value = "Cultura digital II"
$ws.Cells.Item(4, 5).Value = "Cultura digital II"

# --- Column F (Grupo) ---
$ws.Cells.Item(2, 6).Value = "2ARHV"
$ws.Cells.Item(3, 6).Value = "2ARHV"
$ws.Cells.Item(4, 6).Value = "2ALCV"

# --- Column G (Reprobadas) ---
$ws.Cells.Item(2, 7).Value = 4
$ws.Cells.Item(3, 7).Value = 4
$ws.Cells.Item(4, 7).Value = 3
